$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "27.386.05"
$ws.Range("E2").Value = "  +1.29%  "
$ws.Range("D3").Value = "1.861.09"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").Value = "'315.28"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").Value = "'0.4624"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").Value = "'0.3715"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").Value = "'0.07321"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").Value = "'0.8892"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").Value = "'20.12"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").Value = "'0.07820"
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.400"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'6.553"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'91.87"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.761.90"
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("D17").Value = "'1.004"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "'0.000008983"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").Value = "'14.80"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").Value = "27.397.79"
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("D22").Value = "'5.134"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "2.038.93"
$ws.Range("E24").Value = "  +1.93%  "
$ws.Range("D25").Value = "'1.938"
$ws.Range("E25").Value = "  +4.99%  "
$ws.Range("D26").Value = "'152.15"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").Value = "'2.052"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'5.104"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'116.25"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("D31").Value = "'0.08852"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").Value = "'3.137"
$ws.Range("E32").Value = "  +5.32%  "
$ws.Range("D33").Value = "'0.7712"
$ws.Range("E33").Value = "  +5.42%  "
$ws.Range("D34").Value = "'1.173"
$ws.Range("E34").Value = "  +3.39%  "
$ws.Range("D35").Value = "'4.514"
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("D36").Value = "'2.770"
$ws.Range("E36").Value = "  +12.83%  "
$ws.Range("D37").Value = "'1.081"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").Value = "'0.05247"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "'2.970"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").Value = "'7.076"
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("D42").Value = "'0.5146"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("D43").Value = "'0.1640"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").Value = "'8.420"
$ws.Range("E44").Value = "  +2.55%  "
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").Value = "'10.41"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").Value = "'103.26"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").Value = "'1.652"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").Value = "'65.28"
$ws.Range("E51").Value = "  +0.96%  "
